$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A1").Value = "Cat"
$ws.Range("A2").Value = "Persian"
$ws.Range("B1").Value = "Quantity"
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = "Main Coon"
$ws.Range("B3").Value = 3
$ws.Range("A4").Value = "Tabby"
$ws.Range("B4").Value = 1

$ws.Columns.Item(1).ColumnWidth = 12.28515625

$ws.Range("B5").Select()
